$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2669.4375
$ws.Range("I98").Value = 434.75
$ws.Range("J98").Value = 9373.5
$ws.Range("K98").Value = 434.75
$ws.Range("L98").Value = 9373.5
$ws.Range("M98").Value = 1063.25
$ws.Range("N98").Value = -12369.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 7131.125
$ws.Range("I107").Value = 9332.166999999999
$ws.Range("K107").Value = 9332.166999999999
$ws.Range("M107").Value = -7412.166999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2669.4375
$ws.Range("I122").Value = 434.75
$ws.Range("J122").Value = 9373.5
$ws.Range("K122").Value = 1304.25
$ws.Range("L122").Value = 28120.5
$ws.Range("M122").Value = 1145.75
$ws.Range("N122").Value = -33020.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1782.0312
$ws.Range("I137").Value = 896.4761999999999
$ws.Range("J137").Value = 3472.6365
$ws.Range("K137").Value = 2689.4286
$ws.Range("L137").Value = 10417.9095
$ws.Range("M137").Value = -139.4285999999997
$ws.Range("N137").Value = -15517.9095

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2004.9108
$ws.Range("I138").Value = 1146.1471
$ws.Range("J138").Value = 3332.0908
$ws.Range("K138").Value = 3438.4413
$ws.Range("L138").Value = 9996.2724
$ws.Range("M138").Value = 1701.5587
$ws.Range("N138").Value = -20276.2724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3212.318
$ws.Range("I141").Value = 1997.5834
$ws.Range("J141").Value = 4670
$ws.Range("K141").Value = 5992.7502
$ws.Range("L141").Value = 14010
$ws.Range("M141").Value = -812.7502000000004
$ws.Range("N141").Value = -24370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1128
$ws.Range("I2").Value = 1128
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1128
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1015
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8478.989
$ws.Range("I32").Value = 8700.128000000001
$ws.Range("J32").Value = 6750.091
$ws.Range("K32").Value = 8700.128000000001
$ws.Range("L32").Value = 6750.091
$ws.Range("M32").Value = -8413.128000000001
$ws.Range("N32").Value = -7324.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1399.5714
$ws.Range("I45").Value = 1237
$ws.Range("K45").Value = 1237
$ws.Range("M45").Value = -860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I97").Value = 50150
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 50150
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -49654
$ws.Range("N97").Value = -4003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2280
$ws.Range("I102").Value = 2100
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2100
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -478
$ws.Range("N102").Value = -6244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 20000
$ws.Range("J107").Value = 20000
$ws.Range("L107").Value = 20000
$ws.Range("N107").Value = -27680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 33377
$ws.Range("J109").Value = 33377
$ws.Range("L109").Value = 33377
$ws.Range("N109").Value = -36151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1128
$ws.Range("I116").Value = 1128
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1128
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1166
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 18573.143
$ws.Range("I122").Value = 18573.143
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 55719.429
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -53269.429
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10001765
$ws.Range("I132").Value = 13890122
$ws.Range("J132").Value = 3132.7144
$ws.Range("K132").Value = 41670366
$ws.Range("L132").Value = 9398.143199999999
$ws.Range("M132").Value = -41667836
$ws.Range("N132").Value = -14458.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1128
$ws.Range("I3").Value = 1128
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1128
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1014
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 805.1818
$ws.Range("I94").Value = 706.3333
$ws.Range("J94").Value = 1017
$ws.Range("K94").Value = 706.3333
$ws.Range("L94").Value = 1017
$ws.Range("M94").Value = -255.3333
$ws.Range("N94").Value = -1919

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2932.6765
$ws.Range("I58").Value = 679.6957
$ws.Range("J58").Value = 7643.4546
$ws.Range("K58").Value = 679.6957
$ws.Range("L58").Value = 7643.4546
$ws.Range("M58").Value = -476.6957
$ws.Range("N58").Value = -8049.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 28863.666
$ws.Range("J97").Value = 28863.666
$ws.Range("L97").Value = 28863.666
$ws.Range("N97").Value = -30845.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1694.2222
$ws.Range("I99").Value = 1621.1428
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 1621.1428
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = -123.1428000000001
$ws.Range("N99").Value = -4946

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1487.7778
$ws.Range("I105").Value = 1416.6666
$ws.Range("J105").Value = 1630
$ws.Range("K105").Value = 1416.6666
$ws.Range("L105").Value = 1630
$ws.Range("M105").Value = 330.3334
$ws.Range("N105").Value = -5124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1694.2222
$ws.Range("I126").Value = 1621.1428
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 4863.428400000001
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -2393.428400000001
$ws.Range("N126").Value = -10790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2659.9412
$ws.Range("I132").Value = 2555.077
$ws.Range("J132").Value = 3000.75
$ws.Range("K132").Value = 7665.231000000001
$ws.Range("L132").Value = 9002.25
$ws.Range("M132").Value = -5135.231000000001
$ws.Range("N132").Value = -14062.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2932.6765
$ws.Range("I136").Value = 679.6957
$ws.Range("J136").Value = 7643.4546
$ws.Range("K136").Value = 2039.0871
$ws.Range("L136").Value = 22930.3638
$ws.Range("M136").Value = 510.9129
$ws.Range("N136").Value = -28030.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 618.7105
$ws.Range("I113").Value = 497.9
$ws.Range("J113").Value = 752.94446
$ws.Range("K113").Value = 1493.7
$ws.Range("L113").Value = 2258.83338
$ws.Range("M113").Value = 676.3000000000002
$ws.Range("N113").Value = -6598.83338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 53150
$ws.Range("J70").Value = 4200
$ws.Range("L70").Value = 4200
$ws.Range("N70").Value = -4740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 53150
$ws.Range("J73").Value = 4200
$ws.Range("L73").Value = 4200
$ws.Range("N73").Value = -6072

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3510303.2
$ws.Range("I122").Value = 6061650.5
$ws.Range("K122").Value = 18184951.5
$ws.Range("M122").Value = -18182501.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1019.75
$ws.Range("I61").Value = 1019.75
$ws.Range("K61").Value = 1019.75
$ws.Range("M61").Value = -817.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1873.6666
$ws.Range("I100").Value = 1698.4
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 1698.4
$ws.Range("L100").Value = 2750
$ws.Range("M100").Value = -1157.4
$ws.Range("N100").Value = -3832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1019.75
$ws.Range("I113").Value = 1019.75
$ws.Range("K113").Value = 1019.75
$ws.Range("M113").Value = 1150.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6111.5713
$ws.Range("J122").Value = 5384.615
$ws.Range("L122").Value = 16153.845
$ws.Range("N122").Value = -21053.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1871.2727
$ws.Range("I136").Value = 2015
$ws.Range("J136").Value = 1488
$ws.Range("K136").Value = 6045
$ws.Range("L136").Value = 4464
$ws.Range("M136").Value = -3495
$ws.Range("N136").Value = -9564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2225.3704
$ws.Range("I122").Value = 2327.6667
$ws.Range("J122").Value = 2020.7778
$ws.Range("K122").Value = 6983.000100000001
$ws.Range("L122").Value = 6062.3334
$ws.Range("M122").Value = -4533.000100000001
$ws.Range("N122").Value = -10962.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2393.7144
$ws.Range("I126").Value = 1741.4
$ws.Range("J126").Value = 4024.5
$ws.Range("K126").Value = 5224.200000000001
$ws.Range("L126").Value = 12073.5
$ws.Range("M126").Value = -2754.200000000001
$ws.Range("N126").Value = -17013.5
